$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.655.44'
$ws.Range('E2').Value = '  -7.22%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.232.86'
$ws.Range('E3').Value = '  -8.77%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '177.78'
$ws.Range('E5').Value = '  -13.31%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '510.27'
$ws.Range('E6').Value = '  -8.34%  '

$ws.Range('E7').Value = '  -1.86%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.234.86'
$ws.Range('E9').Value = '  -8.47%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.610'
$ws.Range('E10').Value = '  -7.89%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.02'
$ws.Range('E11').Value = '  -7.56%  '

$ws.Range('E12').Value = '  -10.54%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000250'
$ws.Range('E13').Value = '  -8.53%  '

$ws.Range('E14').Value = '  -9.71%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.754.26'
$ws.Range('E15').Value = '  -8.89%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.243.29'
$ws.Range('E16').Value = '  -8.92%  '

$ws.Range('E17').Value = '  -7.50%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '17.44'
$ws.Range('E18').Value = '  -7.11%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '62.598.14'
$ws.Range('E19').Value = '  -7.01%  '

$ws.Range('E20').Value = '  -9.63%  '

$ws.Range('E21').Value = '  -10.76%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '366.13'
$ws.Range('E22').Value = '  -6.79%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.12'
$ws.Range('E23').Value = '  -7.70%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.64'
$ws.Range('E24').Value = '  -10.44%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '78.75'
$ws.Range('E25').Value = '  -5.31%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.02'
$ws.Range('E26').Value = '  -1.99%  '

$ws.Range('E27').Value = '  -0.71%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.23'
$ws.Range('E28').Value = '  -7.39%  '

$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.58'
$ws.Range('E29').Value = '  -8.49%  '

$ws.Range('E30').Value = '  -8.34%  '

$ws.Range('E31').Value = '  -9.28%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '629.08'
$ws.Range('E32').Value = '  -10.40%  '

$ws.Range('E33').Value = '  -9.86%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.05'
$ws.Range('E34').Value = '  -6.83%  '

$ws.Range('E35').Value = '  -6.22%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '58.04'
$ws.Range('E36').Value = '  -8.36%  '

$ws.Range('E37').Value = '  +0.00%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.388'
$ws.Range('E38').Value = '  -5.69%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.70'
$ws.Range('E39').Value = '  -12.21%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.07%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.907.81'
$ws.Range('E41').Value = '  -7.16%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.123'
$ws.Range('E42').Value = '  -5.61%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0₃0639'
$ws.Range('E43').Value = '  -10.54%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.41'
$ws.Range('E44').Value = '  -5.93%  '

$ws.Range('E45').Value = '  -15.07%  '

$ws.Range('E46').Value = '  -5.90%  '

$ws.Range('E47').Value = '  -5.00%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  +4.73%  '

$ws.Range('E49').Value = '  -3.47%  '

$ws.Range('E50').Value = '  -2.80%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.38'
$ws.Range('E51').Value = '  -17.72%  '
